$wb = $excel.ActiveWorkbook

# --- Add "is_targeted list" sheet, placed after "analyte_class list" ---
$afterSheet = $wb.Worksheets.Item("analyte_class list")
$isTargeted = $wb.Worksheets.Add($null, $afterSheet)
$isTargeted.Name = "is_targeted list"
# Use a leading apostrophe so "TRUE"/"FALSE" are stored as text rather than
# being auto-converted to boolean values, then clear the resulting
# quote-prefix formatting so the cells end up with plain default styling.
$isTargeted.Range("A1").Value = "'TRUE"
$isTargeted.Range("A2").Value = "'FALSE"
$isTargeted.Range("A1:A2").ClearFormats()

# --- Add "is_technical_replicate list" sheet, placed after "library_layout list" ---
$afterSheet2 = $wb.Worksheets.Item("library_layout list")
$isTechRep = $wb.Worksheets.Add($null, $afterSheet2)
$isTechRep.Name = "is_technical_replicate list"
$isTechRep.Range("A1").Value = "'TRUE"
$isTechRep.Range("A2").Value = "'FALSE"
$isTechRep.Range("A1:A2").ClearFormats()

# --- Update data validations on the main "Export as TSV" sheet ---
$main = $wb.Worksheets.Item("Export as TSV")

$rngN = $main.Range("N2:N1048576")
$rngN.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$rngN.Validation.ErrorTitle = "Value must come from list"
$rngN.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

$rngAC = $main.Range("AC2:AC1048576")
$rngAC.Validation.Modify(3, 1, 1, "='is_technical_replicate list'!`$A`$1:`$A`$2")
$rngAC.Validation.ErrorTitle = "Value must come from list"
$rngAC.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
